$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("M3_01 Wolf 2021",   13428, 96,  3.723562704795949, 1.489425081918379, 0.4766160262138814),
    @("M3_02 Love 2021",   16016, 96,  3.121878121878122, 1.248751248751249, 0.3996003996003996),
    @("M3_03 Bear 2021",   13636, 192, 3.666764447051921, 1.466705778820768, 0.4693458492226459),
    @("M3_04 Elf 2021",    16049, 96,  3.115458907097016, 1.246183562838806, 0.398778740108418),
    @("M3_05 Viper 2021",  17623, 579, 2.837201384554276, 1.13488055382171,  0.3631617772229473),
    @("M3_06 Magic 2021",  17442, 96,  2.866643733516799, 1.14665749340672,  0.3669303978901502),
    @("M3_07 Griffin 2021",16772, 96,  2.981159074648223, 1.192463629859289, 0.3815883615549726)
)

$startRow = 14
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
}
